$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.931.63"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "2.403.66"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'560.14"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").Value = "'138.31"
$ws.Range("E6").Value = "  +5.27%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").Value = "2.402.91"
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("D11").Value = "'5.71"
$ws.Range("E11").Value = "  +3.90%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").Value = "'25.66"
$ws.Range("E14").Value = "  +7.67%  "
$ws.Range("D15").Value = "2.833.62"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "61.939.76"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").Value = "2.414.47"
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").Value = "'342.53"
$ws.Range("E20").Value = "  +8.49%  "
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("D22").Value = "'6.87"
$ws.Range("E22").Value = "  +3.23%  "
$ws.Range("D24").Value = "'64.83"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'8.28"
$ws.Range("E27").Value = "  +5.22%  "
$ws.Range("D28").Value = "'1.49"
$ws.Range("E28").Value = "  +10.65%  "
$ws.Range("E29").Value = "  +14.43%  "
$ws.Range("D30").Value = "'1.80"
$ws.Range("E30").Value = "  +3.71%  "
$ws.Range("D31").Value = "0.0₃0765"
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").Value = "'6.35"
$ws.Range("E32").Value = "  +6.66%  "
$ws.Range("D33").Value = "'171.67"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "'0.393"
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("E36").Value = "  +11.00%  "
$ws.Range("D37").Value = "'18.48"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'357.76"
$ws.Range("E39").Value = "  +8.73%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'1.65"
$ws.Range("E41").Value = "  +7.45%  "
$ws.Range("D42").Value = "'38.92"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").Value = "'143.55"
$ws.Range("E43").Value = "  +3.60%  "
$ws.Range("D44").Value = "'3.65"
$ws.Range("E44").Value = "  +4.32%  "
$ws.Range("D45").Value = "'20.35"
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("D46").Value = "'0.0960"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").Value = "'0.0516"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("E48").Value = "  +3.81%  "
$ws.Range("D49").Value = "'0.0221"
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("D50").Value = "'17.78"
$ws.Range("E50").Value = "  +5.02%  "
$ws.Range("D51").Value = "0.0₆0216"
$ws.Range("E51").Value = "  -0.82%  "
